$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate a unique piece of text and return the Range of the
# paragraph that contains it (first match only).
# ---------------------------------------------------------------------------
function Get-ParagraphContaining($anchorText) {
    $r = $d.Content
    $found = $r.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $anchorText"
    }
    return $r.Paragraphs(1).Range
}

# ---------------------------------------------------------------------------
# Helper: force the engine to coalesce a span of (identically formatted)
# runs into a single run without altering its visible text. Setting a
# Range's .Text only rewrites run boundaries when the assigned text
# actually differs from what is already there, so the helper briefly
# perturbs the text by one character and then restores the original.
# ---------------------------------------------------------------------------
function Merge-Range($rng) {
    $original = $rng.Text
    $rng.Text = $original + "X"
    $tmp = $d.Range($rng.Start, $rng.Start + $original.Length + 1)
    $tmp.Text = $original
}

# ===========================================================================
# Paragraph: "Il software viene eseguito su EC2 ..."
# Net visible change: a comma is inserted after "Amazon AWS" (right before
# " ed è uno degli script ..."). The containing runs are also re-merged
# while making the edit, matching what Word itself does on save.
# ===========================================================================
$p5Range = Get-ParagraphContaining "Il software viene eseguito su EC2"
$p5Start = $p5Range.Start
$p5TextRaw = $p5Range.Text
$p5Text = $p5TextRaw.Substring(0, $p5TextRaw.Length - 1)   # drop trailing paragraph mark

$idxATutti = $p5Text.IndexOf(" a tutti gli utenti")
$idxTelegram = $p5Text.IndexOf("Telegram")
$idxEc2End = $p5Text.IndexOf(", gestit")

# 1) Merge the tail: " a tutti gli utenti ... È in ascolto sulla porta 5001."
$tailStart = $p5Start + $idxATutti
$tailEnd = $p5Start + $p5Text.Length
Merge-Range ($d.Range($tailStart, $tailEnd))

# 2) Rewrite the middle block that currently reads
#    ", gestita dal sistema di Amazon AWS ed è uno degli script ... tramite "
#    into ", ed è uno degli script ... tramite " (the leading "gestita dal
#    sistema di Amazon AWS" becomes redundant once it's duplicated onto the
#    end of the opening sentence in step 3 below).
$midStart = $p5Start + $idxEc2End + 1   # +1 to skip the existing ","
$midEnd = $p5Start + $idxTelegram
$midRange = $d.Range($midStart, $midEnd)
$midRange.Text = " ed è uno degli script che compongono la parte Cloud del progetto. È collegato al database MySQL presente sulla stessa macchina e riceve ed invia messaggi tramite "

# 3) Extend the opening sentence so it already announces the AWS/EC2 detail:
#    "Il software viene eseguito su EC2" -> "Il software viene eseguito su EC2, gestita dal sistema di Amazon AWS"
$openRange = $d.Range($p5Start, $p5Start + $idxEc2End)
$openRange.Text = "Il software viene eseguito su EC2, gestita dal sistema di Amazon AWS"

# ===========================================================================
# "pymysql" bullet: merge the trailing runs (no visible text change)
#   ", per comunicare con il database" + ", " + "reperire le informazioni" + " ed aggiornarlo."
# ===========================================================================
$pyParaRange = Get-ParagraphContaining "pymysql"
$pyParaText = $pyParaRange.Text.Substring(0, $pyParaRange.Text.Length - 1)
$pyAfter = $pyParaRange.Start + $pyParaText.IndexOf("pymysql") + "pymysql".Length
$pyEnd = $pyParaRange.Start + $pyParaText.Length
Merge-Range ($d.Range($pyAfter, $pyEnd))

# ===========================================================================
# "requests" bullet: merge the runs (no visible text change)
#   ", per" + " ricevere" + " la richiesta POST " + "dello script “"
# ===========================================================================
$reqParaRange = Get-ParagraphContaining "requests"
$reqParaText = $reqParaRange.Text.Substring(0, $reqParaRange.Text.Length - 1)
$reqAfter = $reqParaRange.Start + $reqParaText.IndexOf("requests") + "requests".Length
$reqDelloEnd = $reqParaRange.Start + $reqParaText.IndexOf("dello script “") + "dello script “".Length
Merge-Range ($d.Range($reqAfter, $reqDelloEnd))

Write-Host "Edits applied."
